$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Addr, $Val)
    $cell = $ws.Range($Addr)
    $cell.NumberFormat = "@"
    $cell.Value = $Val
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '35.293.14'
Set-TextValue 'E2' '  -0.73%  '
Set-TextValue 'D3' '1.898.24'
Set-TextValue 'E3' '  -0.35%  '
Set-TextValue 'E4' '  -0.23%  '
Set-TextValue 'D5' '246.91'
Set-TextValue 'E5' '  -0.25%  '
Set-TextValue 'D6' '0.694'
Set-TextValue 'E6' '  +9.65%  '
Set-TextValue 'D8' '40.56'
Set-TextValue 'E8' '  -3.85%  '
Set-TextValue 'D9' '0.349'
Set-TextValue 'E9' '  +2.83%  '
Set-TextValue 'D10' '52.33'
Set-TextValue 'E10' '  +8.17%  '
Set-TextValue 'E11' '  +2.44%  '
Set-TextValue 'D12' '0.0988'
Set-TextValue 'E12' '  -1.03%  '
Set-TextValue 'D13' '2.172.43'
Set-TextValue 'E13' '  -0.38%  '
Set-TextValue 'D14' '12.59'
Set-TextValue 'E14' '  +1.51%  '
Set-TextValue 'D15' '0.709'
Set-TextValue 'E15' '  +2.64%  '
Set-TextValue 'D16' '1.894.42'
Set-TextValue 'E16' '  -0.63%  '
Set-TextValue 'D17' '4.83'
Set-TextValue 'E17' '  -0.62%  '
Set-TextValue 'D18' '35.266.75'
Set-TextValue 'E18' '  -0.72%  '
Set-TextValue 'D19' '72.27'
Set-TextValue 'E19' '  +0.26%  '
Set-TextValue 'D20' '0.0₃0821'
Set-TextValue 'E20' '  +0.74%  '
Set-TextValue 'D21' '240.97'
Set-TextValue 'E21' '  -1.27%  '
Set-TextValue 'D22' '12.75'
Set-TextValue 'E22' '  +2.05%  '
Set-TextValue 'D23' '4.80'
Set-TextValue 'E23' '  -1.80%  '
Set-TextValue 'E24' '  -0.20%  '
Set-TextValue 'E25' '  +1.85%  '
Set-TextValue 'D26' '2.33'
Set-TextValue 'E26' '  +3.20%  '
Set-TextValue 'D27' '168.79'
Set-TextValue 'E27' '  -1.89%  '
Set-TextValue 'D28' '8.65'
Set-TextValue 'E28' '  +0.54%  '
Set-TextValue 'D29' '19.15'
Set-TextValue 'E29' '  +6.48%  '
Set-TextValue 'D30' '0.131'
Set-TextValue 'E30' '  +4.48%  '
Set-TextValue 'E31' '  +20.05%  '
Set-TextValue 'D32' '4.17'
Set-TextValue 'E32' '  +1.48%  '
Set-TextValue 'D33' '0.0568'
Set-TextValue 'E33' '  +0.23%  '
Set-TextValue 'B34' 'BinanceUSD'
Set-TextValue 'C34' 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue 'D34' '1.01'
Set-TextValue 'E34' '  -0.30%  '
Set-TextValue 'B35' 'WEMIXToken'
Set-TextValue 'C35' 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D35' '1.87'
Set-TextValue 'E35' '  +7.21%  '
Set-TextValue 'D36' '4.12'
Set-TextValue 'E36' '  -1.66%  '
Set-TextValue 'B37' 'ImmutableX'
Set-TextValue 'C37' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D37' '0.911'
Set-TextValue 'E37' '  -6.25%  '
Set-TextValue 'B38' 'TrustWalletToken'
Set-TextValue 'C38' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D38' '1.50'
Set-TextValue 'E38' '  +12.96%  '
Set-TextValue 'E39' '  +0.23%  '
Set-TextValue 'D40' '0.0658'
Set-TextValue 'E40' '  +11.09%  '
Set-TextValue 'B41' 'ARBITRUM'
Set-TextValue 'C41' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D41' '1.09'
Set-TextValue 'E41' '  -1.13%  '
Set-TextValue 'B42' 'InjectiveProtocol'
Set-TextValue 'C42' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D42' '16.45'
Set-TextValue 'E42' '  +5.66%  '
Set-TextValue 'D43' '93.57'
Set-TextValue 'E43' '  +2.73%  '
Set-TextValue 'E44' '  +0.79%  '
Set-TextValue 'D45' '1.350.66'
Set-TextValue 'E45' '  -0.38%  '
Set-TextValue 'E46' '  +2.62%  '
Set-TextValue 'E47' '  +0.59%  '
Set-TextValue 'E48' '  +0.90%  '
Set-TextValue 'D49' '12.28'
Set-TextValue 'E49' '  -2.40%  '
Set-TextValue 'D50' '45.04'
Set-TextValue 'E50' '  -7.91%  '
Set-TextValue 'D51' '6.47'
